# Auto-generated Excel COM-interop script
# Refreshes the crypto price/volume table: updates Price (D) and Volume(1h) (E)
# values for most rows, and reorders three coin rows (12-14: WrappedEther /
# WrappedliquidstakedEther2.0 / Polkadot) and swaps two rows (20-21: BitcoinCash /
# Uniswap) including their Coin name (B) and Link (C) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.846.66'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.631.59'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  +0.45%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("E6").Value = '  +0.97%  '
$ws.Range("E7").Value = '  +0.51%  '
$ws.Range("E8").Value = '  -0.51%  '
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("E10").Value = '  -0.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0791'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.675.89'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.16%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.856.70'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.24'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("E15").Value = '  -1.16%  '
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.51'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.859.62'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("E19").Value = '  +0.46%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.54%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.79%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("E24").Value = '  +1.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.95%  '
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("E27").Value = '  +2.36%  '
$ws.Range("E29").Value = '  -0.57%  '
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("E31").Value = '  +1.25%  '
$ws.Range("E32").Value = '  -0.52%  '
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("E34").Value = '  -2.12%  '
$ws.Range("E35").Value = '  +1.57%  '
$ws.Range("E36").Value = '  -0.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.138.14'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.29%  '
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("E40").Value = '  +0.52%  '
$ws.Range("E41").Value = '  +0.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.26'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.42'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.797'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.765.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.19'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.55%  '
$ws.Range("E48").Value = '  +3.24%  '
$ws.Range("E49").Value = '  -1.11%  '
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.62'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.99%  '
